$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the default (unstyled) style from a cell that keeps default formatting,
# so we can force text values via an apostrophe prefix without permanently altering cell style.
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Formula = "'26.555.86"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Formula = "'  -2.45%  "
$ws.Range("E2").Style = $defaultStyle
$ws.Range("D3").Formula = "'1.815.56"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Formula = "'  -2.07%  "
$ws.Range("E3").Style = $defaultStyle
$ws.Range("D4").Formula = "'1.009"
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").Formula = "'  +0.79%  "
$ws.Range("E4").Style = $defaultStyle
$ws.Range("D5").Formula = "'308.71"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Formula = "'  -1.65%  "
$ws.Range("E5").Style = $defaultStyle
$ws.Range("E6").Formula = "'  +0.62%  "
$ws.Range("E6").Style = $defaultStyle
$ws.Range("E7").Formula = "'  -1.46%  "
$ws.Range("E7").Style = $defaultStyle
$ws.Range("D8").Formula = "'0.3667"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Formula = "'  -1.26%  "
$ws.Range("E8").Style = $defaultStyle
$ws.Range("D9").Formula = "'0.07142"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Formula = "'  -2.04%  "
$ws.Range("E9").Style = $defaultStyle
$ws.Range("D10").Formula = "'0.8776"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Formula = "'  -1.03%  "
$ws.Range("E10").Style = $defaultStyle
$ws.Range("D11").Formula = "'0.07776"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Formula = "'  -1.07%  "
$ws.Range("E11").Style = $defaultStyle
$ws.Range("D12").Formula = "'19.36"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Formula = "'  -3.64%  "
$ws.Range("E12").Style = $defaultStyle
$ws.Range("D13").Formula = "'1.804.53"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Formula = "'  -3.24%  "
$ws.Range("E13").Style = $defaultStyle
$ws.Range("D14").Formula = "'5.293"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Formula = "'  -1.88%  "
$ws.Range("E14").Style = $defaultStyle
$ws.Range("D15").Formula = "'6.374"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("D16").Formula = "'86.30"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Formula = "'  -5.23%  "
$ws.Range("E16").Style = $defaultStyle
$ws.Range("D17").Formula = "'1.010"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Formula = "'  +0.84%  "
$ws.Range("E17").Style = $defaultStyle
$ws.Range("D18").Formula = "'0.000008618"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E19").Formula = "'  +0.66%  "
$ws.Range("E19").Style = $defaultStyle
$ws.Range("D20").Formula = "'26.645.70"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Formula = "'  -2.24%  "
$ws.Range("E20").Style = $defaultStyle
$ws.Range("D21").Formula = "'14.27"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Formula = "'  -2.94%  "
$ws.Range("E21").Style = $defaultStyle
$ws.Range("D22").Formula = "'5.005"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Formula = "'  -1.59%  "
$ws.Range("E22").Style = $defaultStyle
$ws.Range("E23").Formula = "'  -0.44%  "
$ws.Range("E23").Style = $defaultStyle
$ws.Range("D24").Formula = "'1.989"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Formula = "'  +1.95%  "
$ws.Range("E24").Style = $defaultStyle
$ws.Range("D25").Formula = "'151.61"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Formula = "'  +0.23%  "
$ws.Range("E25").Style = $defaultStyle
$ws.Range("D26").Formula = "'17.97"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Formula = "'  -2.30%  "
$ws.Range("E26").Style = $defaultStyle
$ws.Range("D27").Formula = "'2.077"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Formula = "'  +1.42%  "
$ws.Range("E27").Style = $defaultStyle
$ws.Range("D28").Formula = "'113.11"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Formula = "'  -2.40%  "
$ws.Range("E28").Style = $defaultStyle
$ws.Range("D29").Formula = "'4.862"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Formula = "'  -3.61%  "
$ws.Range("E29").Style = $defaultStyle
$ws.Range("D30").Formula = "'0.08696"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Formula = "'  -1.19%  "
$ws.Range("E30").Style = $defaultStyle
$ws.Range("D31").Formula = "'3.062"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Formula = "'  -2.49%  "
$ws.Range("E31").Style = $defaultStyle
$ws.Range("D32").Formula = "'4.514"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Formula = "'  -0.09%  "
$ws.Range("E32").Style = $defaultStyle
$ws.Range("D33").Formula = "'0.7352"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Formula = "'  -4.34%  "
$ws.Range("E33").Style = $defaultStyle
$ws.Range("D34").Formula = "'2.694"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Formula = "'  -1.40%  "
$ws.Range("E34").Style = $defaultStyle
$ws.Range("D35").Formula = "'1.120"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Formula = "'  -3.98%  "
$ws.Range("E35").Style = $defaultStyle
$ws.Range("E36").Formula = "'  +0.50%  "
$ws.Range("E36").Style = $defaultStyle
$ws.Range("D37").Formula = "'1.085"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Formula = "'  -1.68%  "
$ws.Range("E37").Style = $defaultStyle
$ws.Range("D38").Formula = "'0.01946"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Formula = "'  +0.46%  "
$ws.Range("E38").Style = $defaultStyle
$ws.Range("D39").Formula = "'0.05123"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Formula = "'  -1.87%  "
$ws.Range("E39").Style = $defaultStyle
$ws.Range("D40").Formula = "'2.910"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Formula = "'  -1.00%  "
$ws.Range("E40").Style = $defaultStyle
$ws.Range("D41").Formula = "'7.010"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Formula = "'  -0.37%  "
$ws.Range("E41").Style = $defaultStyle
$ws.Range("D42").Formula = "'0.5033"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Formula = "'  -1.75%  "
$ws.Range("E42").Style = $defaultStyle
$ws.Range("D43").Formula = "'0.1558"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Formula = "'  -4.26%  "
$ws.Range("E43").Style = $defaultStyle
$ws.Range("D44").Formula = "'8.186"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Formula = "'  -3.11%  "
$ws.Range("E44").Style = $defaultStyle
$ws.Range("D45").Formula = "'1.007"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Formula = "'  +0.68%  "
$ws.Range("E45").Style = $defaultStyle
$ws.Range("D46").Formula = "'0.4628"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Formula = "'  -3.54%  "
$ws.Range("E46").Style = $defaultStyle
$ws.Range("E47").Formula = "'  -3.45%  "
$ws.Range("E47").Style = $defaultStyle
$ws.Range("D48").Formula = "'101.34"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Formula = "'  -1.49%  "
$ws.Range("E48").Style = $defaultStyle
$ws.Range("D49").Formula = "'1.596"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Formula = "'  -2.86%  "
$ws.Range("E49").Style = $defaultStyle
$ws.Range("D50").Formula = "'0.06007"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Formula = "'  -3.22%  "
$ws.Range("E50").Style = $defaultStyle
$ws.Range("D51").Formula = "'64.53"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Formula = "'  -1.51%  "
$ws.Range("E51").Style = $defaultStyle
